$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data for column E (NEC method) and column F (LBNL method) ---
# Row 3
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 5

# Row 4
$ws.Range("E4").Value = 3.15
$ws.Range("F4").Value = 6.61

# Row 5 (text "Error" entries, same as column D)
$ws.Range("E5").Value = "Error"
$ws.Range("F5").Value = "Error"

# Row 6
$ws.Range("E6").Value = 3.4
$ws.Range("F6").Value = 6.61

# Row 7
$ws.Range("E7").Value = 4.21
$ws.Range("F7").NumberFormat = "0.00"
$ws.Range("F7").Value = 6.7220000000000004

# Row 8 (E8 already formatted with 0.00; add F8 matching style)
$ws.Range("F8").NumberFormat = "0.00"
$ws.Range("F8").Value = 8.8733799999999992

# Row 9
$ws.Range("F9").Value = 13.44

# Row 10
$ws.Range("E10").Value = 9.4600000000000009
$ws.Range("F10").Value = 9.4600000000000009

# Row 11
$ws.Range("E11").Value = 27.85
$ws.Range("F11").Value = 27.85

# Row 12 (E12 value changes from 1.9869000000000001 to 7.95, F12 is new)
$ws.Range("E12").Value = 7.95
$ws.Range("F12").NumberFormat = "0.00"
$ws.Range("F12").Value = 7.95

# --- Update the selected cell in the sheet view ---
$ws.Range("K10").Select()
